$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 453.3125
$ws.Range("J33").Value = 701.7143
$ws.Range("L33").Value = 701.7143
$ws.Range("N33").Value = -1159.7143
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("H70").Value = 4856.7144
$ws.Range("J70").Value = 4856.7144
$ws.Range("L70").Value = 14570.1432
$ws.Range("N70").Value = -15110.1432
$ws.Range("H73").Value = 4856.7144
$ws.Range("J73").Value = 4856.7144
$ws.Range("L73").Value = 14570.1432
$ws.Range("N73").Value = -16442.1432
$ws.Range("H116").Value = 5749
$ws.Range("I116").Value = 5748
$ws.Range("J116").Value = 5750
$ws.Range("K116").Value = 5748
$ws.Range("L116").Value = 5750
$ws.Range("M116").Value = -2306
$ws.Range("N116").Value = -12634
$ws.Range("H137").Value = 1770.96
$ws.Range("I137").Value = 1363.2354
$ws.Range("K137").Value = 4089.7062
$ws.Range("M137").Value = -1539.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 1100
$ws.Range("I30").Value = 1100
$ws.Range("K30").Value = 1100
$ws.Range("M30").Value = -950
$ws.Range("H110").Value = 3416.0908
$ws.Range("I110").Value = 644.25
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 644.25
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = 1400.75
$ws.Range("N110").Value = -9090
$ws.Range("H119").Value = 52000
$ws.Range("J119").Value = 52000
$ws.Range("L119").Value = 52000
$ws.Range("N119").Value = -61676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27221.4
$ws.Range("I82").Value = 8053.5
$ws.Range("K82").Value = 8053.5
$ws.Range("M82").Value = -7670.5
$ws.Range("H85").Value = 27221.4
$ws.Range("I85").Value = 8053.5
$ws.Range("K85").Value = 8053.5
$ws.Range("M85").Value = -6727.5
$ws.Range("H105").Value = 704
$ws.Range("I105").Value = 700
$ws.Range("J105").Value = 708
$ws.Range("K105").Value = 700
$ws.Range("L105").Value = 708
$ws.Range("M105").Value = 1047
$ws.Range("N105").Value = -4202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 825
$ws.Range("I22").Value = 900
$ws.Range("K22").Value = 900
$ws.Range("M22").Value = -550
$ws.Range("H43").Value = 28666.666
$ws.Range("J43").Value = 28666.666
$ws.Range("L43").Value = 28666.666
$ws.Range("N43").Value = -29034.666
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
$ws.Range("H59").Value = 34571.285
$ws.Range("I59").Value = 32000
$ws.Range("J59").Value = 34999.832
$ws.Range("K59").Value = 32000
$ws.Range("L59").Value = 34999.832
$ws.Range("M59").Value = -30855
$ws.Range("N59").Value = -37289.832
$ws.Range("H68").Value = 37533.5
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498
$ws.Range("H71").Value = 37533.5
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488
$ws.Range("H74").Value = 39499.145
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748
$ws.Range("H77").Value = 39499.145
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736
$ws.Range("H101").Value = 28666.666
$ws.Range("J101").Value = 28666.666
$ws.Range("L101").Value = 28666.666
$ws.Range("N101").Value = -35156.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1202.3556
$ws.Range("I4").Value = 1084.2069
$ws.Range("J4").Value = 1416.5
$ws.Range("K4").Value = 3252.620699999999
$ws.Range("L4").Value = 4249.5
$ws.Range("M4").Value = -3140.620699999999
$ws.Range("N4").Value = -4473.5
$ws.Range("H113").Value = 648.1667
$ws.Range("I113").Value = 345
$ws.Range("K113").Value = 1035
$ws.Range("M113").Value = 1135

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12500
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15302
$ws.Range("H80").Value = 2807.8
$ws.Range("I80").Value = 995
$ws.Range("J80").Value = 3009.2222
$ws.Range("K80").Value = 995
$ws.Range("L80").Value = 3009.2222
$ws.Range("M80").Value = 3
$ws.Range("N80").Value = -5005.2222
$ws.Range("H83").Value = 2807.8
$ws.Range("I83").Value = 995
$ws.Range("J83").Value = 3009.2222
$ws.Range("K83").Value = 4975
$ws.Range("L83").Value = 15046.111
$ws.Range("M83").Value = 17
$ws.Range("N83").Value = -25030.111
$ws.Range("H97").Value = 1067.9
$ws.Range("I97").Value = 878.1667
$ws.Range("K97").Value = 878.1667
$ws.Range("M97").Value = -382.1667
$ws.Range("H102").Value = 2332.6667
$ws.Range("J102").Value = 2999
$ws.Range("L102").Value = 2999
$ws.Range("N102").Value = -6243

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1894
$ws.Range("I22").Value = 1726
$ws.Range("J22").Value = 2650
$ws.Range("K22").Value = 1726
$ws.Range("L22").Value = 2650
$ws.Range("M22").Value = -1431
$ws.Range("N22").Value = -3240
$ws.Range("H26").Value = 8000
$ws.Range("I26").Value = 8000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -7705
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 1894
$ws.Range("I27").Value = 1726
$ws.Range("J27").Value = 2650
$ws.Range("K27").Value = 1726
$ws.Range("L27").Value = 2650
$ws.Range("M27").Value = -1619
$ws.Range("N27").Value = -2864
$ws.Range("H31").Value = 1517.8334
$ws.Range("I31").Value = 1221.4
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1221.4
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -973.4000000000001
$ws.Range("N31").Value = -3496
$ws.Range("H55").Value = 231.2
$ws.Range("I55").Value = 151.25
$ws.Range("K55").Value = 151.25
$ws.Range("M55").Value = 21.75
$ws.Range("H94").Value = 90000
$ws.Range("J94").Value = 90000
$ws.Range("L94").Value = 90000
$ws.Range("N94").Value = -91352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 21360.8
$ws.Range("J113").Value = 1750.5
$ws.Range("L113").Value = 5251.5
$ws.Range("N113").Value = -9591.5

Write-Host "Applied Marilith_Profits updates"